$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Cells whose new value is a plain (non-numeric-looking) string literal -
# Excel keeps these as text automatically.
$ws.Range('D2').Value = '43.747.11'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.320.16'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  +5.60%  '
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.668.00'
$ws.Range('E14').Value = '  +4.05%  '
$ws.Range('E16').Value = '  +8.19%  '
$ws.Range('D17').Value = '2.327.82'
$ws.Range('E17').Value = '  +4.50%  '
$ws.Range('D18').Value = '43.779.41'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('E19').Value = '  +5.96%  '
$ws.Range('E20').Value = '  +7.17%  '
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('E22').Value = '  +2.71%  '
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('E24').Value = '  +3.95%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('E28').Value = '  -2.29%  '
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('E30').Value = '  -5.36%  '
$ws.Range('E31').Value = '  +7.74%  '
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  +3.17%  '
$ws.Range('E36').Value = '  +2.08%  '
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  -5.39%  '
$ws.Range('E40').Value = '  +11.28%  '
$ws.Range('E41').Value = '  +8.47%  '
$ws.Range('E42').Value = '  +19.14%  '
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('E44').Value = '  +9.80%  '
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('E47').Value = '  +4.08%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('E50').Value = '  +16.96%  '
$ws.Range('D51').Value = '2.546.79'
$ws.Range('E51').Value = '  +3.83%  '

# Cells whose new value looks like a valid number (e.g. "97.39") - these
# must be forced to Text format first, or Excel auto-converts them to
# a numeric cell and the formatting (trailing zero, exact digits) is lost.
Set-TextCell 'D5' '97.39'
Set-TextCell 'D6' '270.88'
Set-TextCell 'D9' '0.623'
Set-TextCell 'D10' '45.76'
Set-TextCell 'D11' '0.0948'
Set-TextCell 'D12' '8.11'
Set-TextCell 'D15' '15.49'
Set-TextCell 'D21' '72.73'
Set-TextCell 'D22' '239.23'
Set-TextCell 'D24' '9.42'
Set-TextCell 'D27' '11.32'
Set-TextCell 'D28' '3.46'
Set-TextCell 'D30' '38.26'
Set-TextCell 'D31' '22.40'
Set-TextCell 'D32' '174.96'
Set-TextCell 'D33' '0.0904'
Set-TextCell 'D34' '5.47'
Set-TextCell 'D36' '0.0359'
Set-TextCell 'D37' '0.108'
Set-TextCell 'D38' '4.40'
Set-TextCell 'D40' '0.243'
Set-TextCell 'D42' '1.37'
Set-TextCell 'D43' '12.18'
Set-TextCell 'D44' '9.19'
Set-TextCell 'D45' '62.04'
Set-TextCell 'D46' '5.37'
Set-TextCell 'D48' '100.40'
